$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the 5 kept product rows (8-12) ---
# Row 8
$ws.Range("B8").Value = 5
$ws.Range("D8").Value = 42
$ws.Range("F8").Value = "GENERAL 11KW 0 KONTAKTÖR"
$ws.Range("G8").Value = "Adet"
$ws.Range("H8").Value = 100

# Row 9
$ws.Range("B9").Value = 5
$ws.Range("D9").Value = 43
$ws.Range("F9").Value = "4X16+10 ALÜMİNYUM KABLO"
$ws.Range("G9").Value = "Adet"
$ws.Range("H9").Value = 200

# Row 10
$ws.Range("B10").Value = 5
$ws.Range("D10").Value = 44
$ws.Range("F10").Value = "DILM225-S/22 KONTAKTÖR 110KW"
$ws.Range("G10").Value = "Adet"
$ws.Range("H10").Value = 300

# Row 11
$ws.Range("B11").Value = 5
$ws.Range("D11").Value = 45
$ws.Range("F11").Value = "C63-3X63A GRUP OTOMAT"
$ws.Range("G11").Value = "Adet"
$ws.Range("H11").Value = 400

# Row 12
$ws.Range("B12").Value = 5
$ws.Range("D12").Value = 46
$ws.Range("F12").Value = "4X1,5 NYM KABLO"
$ws.Range("G12").Value = "Adet"
$ws.Range("H12").Value = 500

# --- Clear out the now-unused product rows (13-38) ---
# B:D (hidden helper columns) are removed entirely; E:H are cleared but
# keep their row styling.
$ws.Range("B13:D38").ClearContents()
$ws.Range("E13:H38").ClearContents()
